$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.681.68"
$ws.Range("E2").Value = "  -3.51%  "
$ws.Range("D3").Value = "1.738.97"
$ws.Range("D5").Value = "'238.13"
$ws.Range("E5").Value = "  -8.19%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  -7.06%  "
$ws.Range("D8").Value = "'41.87"
$ws.Range("E8").Value = "  -7.11%  "
$ws.Range("D9").Value = "'0.2417"
$ws.Range("E9").Value = "  -23.38%  "
$ws.Range("D10").Value = "'0.05984"
$ws.Range("E10").Value = "  -12.10%  "
$ws.Range("D11").Value = "1.736.24"
$ws.Range("E11").Value = "  -5.72%  "
$ws.Range("D12").Value = "'0.06746"
$ws.Range("E12").Value = "  -13.02%  "
$ws.Range("E13").Value = "  -21.74%  "
$ws.Range("E14").Value = "  -11.59%  "
$ws.Range("D15").Value = "'0.5824"
$ws.Range("E15").Value = "  -25.46%  "
$ws.Range("D16").Value = "'76.17"
$ws.Range("E16").Value = "  -13.48%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "25.712.44"
$ws.Range("E19").Value = "  -3.51%  "
$ws.Range("D20").Value = "'11.44"
$ws.Range("E20").Value = "  -17.42%  "
$ws.Range("D21").Value = "'0.000006333"
$ws.Range("E21").Value = "  -20.21%  "
$ws.Range("D22").Value = "1.955.56"
$ws.Range("E22").Value = "  -5.71%  "
$ws.Range("D23").Value = "'3.935"
$ws.Range("E23").Value = "  -14.67%  "
$ws.Range("D24").Value = "'5.101"
$ws.Range("E24").Value = "  -14.73%  "
$ws.Range("D25").Value = "'7.822"
$ws.Range("E25").Value = "  -16.04%  "
$ws.Range("D26").Value = "'135.80"
$ws.Range("E26").Value = "  -4.56%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'1.456"
$ws.Range("E27").Value = "  -13.14%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'1.831"
$ws.Range("E28").Value = "  -17.36%  "
$ws.Range("D29").Value = "'14.39"
$ws.Range("E29").Value = "  -15.34%  "
$ws.Range("D30").Value = "'99.88"
$ws.Range("E30").Value = "  -9.91%  "
$ws.Range("D31").Value = "'0.08099"
$ws.Range("E31").Value = "  -7.27%  "
$ws.Range("D32").Value = "'3.698"
$ws.Range("E32").Value = "  -11.74%  "
$ws.Range("D33").Value = "'3.357"
$ws.Range("E33").Value = "  -17.66%  "
$ws.Range("D34").Value = "'0.04335"
$ws.Range("E34").Value = "  -11.08%  "
$ws.Range("D35").Value = "'1.001"
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").Value = "'2.664"
$ws.Range("E36").Value = "  -6.84%  "
$ws.Range("D37").Value = "'1.020"
$ws.Range("E37").Value = "  -10.56%  "
$ws.Range("D38").Value = "'0.5983"
$ws.Range("E38").Value = "  -18.07%  "
$ws.Range("D39").Value = "'2.731"
$ws.Range("E39").Value = "  -11.64%  "
$ws.Range("D40").Value = "'2.025"
$ws.Range("E40").Value = "  -10.41%  "
$ws.Range("D41").Value = "'1.001"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "'102.66"
$ws.Range("E42").Value = "  -6.65%  "
$ws.Range("E43").Value = "  -14.54%  "
$ws.Range("D44").Value = "'0.7931"
$ws.Range("E44").Value = "  -11.53%  "
$ws.Range("D45").Value = "'0.3789"
$ws.Range("E45").Value = "  -20.86%  "
$ws.Range("D46").Value = "'5.104"
$ws.Range("E46").Value = "  -13.92%  "
$ws.Range("D47").Value = "'5.999"
$ws.Range("E47").Value = "  -21.84%  "
$ws.Range("D48").Value = "'0.05090"
$ws.Range("D49").Value = "'30.02"
$ws.Range("E49").Value = "  -13.70%  "
$ws.Range("D50").Value = "'0.1054"
$ws.Range("E50").Value = "  -14.96%  "
$ws.Range("D51").Value = "'52.07"
$ws.Range("E51").Value = "  -12.86%  "
